$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 1268.9
$ws.Range("J112").Value = 1312.3243
$ws.Range("L112").Value = 3936.9729
$ws.Range("N112").Value = -6152.9729
# Row 132
$ws.Range("H132").Value = 3207.6965
$ws.Range("I132").Value = 2364.62
$ws.Range("J132").Value = 10233.333
$ws.Range("K132").Value = 7093.86
$ws.Range("L132").Value = 30699.999
$ws.Range("M132").Value = -4563.86
$ws.Range("N132").Value = -35759.999
# Row 138
$ws.Range("H138").Value = 2533.7
$ws.Range("I138").Value = 926.6539
$ws.Range("J138").Value = 3098.338
$ws.Range("K138").Value = 2779.9617
$ws.Range("L138").Value = 9295.014000000001
$ws.Range("M138").Value = 2360.0383
$ws.Range("N138").Value = -19575.014

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7915.36
$ws.Range("I32").Value = 6421.161
$ws.Range("J32").Value = 17915
$ws.Range("K32").Value = 6421.161
$ws.Range("L32").Value = 17915
$ws.Range("M32").Value = -6134.161
$ws.Range("N32").Value = -18489
# Row 43
$ws.Range("H43").Value = 4966.6665
$ws.Range("J43").Value = 4966.6665
$ws.Range("L43").Value = 4966.6665
$ws.Range("N43").Value = -5592.6665

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 570887.9
$ws.Range("I31").Value = 2008.2
$ws.Range("J31").Value = 1217342
$ws.Range("K31").Value = 2008.2
$ws.Range("L31").Value = 1217342
$ws.Range("M31").Value = -1713.2
$ws.Range("N31").Value = -1217932
# Row 34
$ws.Range("H34").Value = 570887.9
$ws.Range("I34").Value = 2008.2
$ws.Range("J34").Value = 1217342
$ws.Range("K34").Value = 2008.2
$ws.Range("L34").Value = 1217342
$ws.Range("M34").Value = -1806.2
$ws.Range("N34").Value = -1217746
# Row 58
$ws.Range("H58").Value = 1726.4419
$ws.Range("I58").Value = 1145.9395
$ws.Range("J58").Value = 3642.1
$ws.Range("K58").Value = 1145.9395
$ws.Range("L58").Value = 3642.1
$ws.Range("M58").Value = -942.9395
$ws.Range("N58").Value = -4048.1
# Row 99
$ws.Range("H99").Value = 1248.1428
$ws.Range("I99").Value = 1190.3077
$ws.Range("K99").Value = 1190.3077
$ws.Range("M99").Value = 307.6922999999999
# Row 126
$ws.Range("H126").Value = 1248.1428
$ws.Range("I126").Value = 1190.3077
$ws.Range("K126").Value = 3570.9231
$ws.Range("M126").Value = -1100.9231
# Row 134
$ws.Range("H134").Value = 3478.05
$ws.Range("I134").Value = 1611.6666
$ws.Range("J134").Value = 5005.091
$ws.Range("K134").Value = 4834.9998
$ws.Range("L134").Value = 15015.273
$ws.Range("M134").Value = -2299.9998
$ws.Range("N134").Value = -20085.273
# Row 136
$ws.Range("H136").Value = 1726.4419
$ws.Range("I136").Value = 1145.9395
$ws.Range("J136").Value = 3642.1
$ws.Range("K136").Value = 3437.8185
$ws.Range("L136").Value = 10926.3
$ws.Range("M136").Value = -887.8184999999999
$ws.Range("N136").Value = -16026.3

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 2203.1052
$ws.Range("I39").Value = 700
$ws.Range("J39").Value = 2379.9412
$ws.Range("K39").Value = 2100
$ws.Range("L39").Value = 7139.823600000001
$ws.Range("M39").Value = -1806
$ws.Range("N39").Value = -7727.823600000001
# Row 107
$ws.Range("H107").Value = 625394.25
$ws.Range("I107").Value = 148.375
$ws.Range("K107").Value = 445.125
$ws.Range("M107").Value = 1474.875
# Row 113
$ws.Range("H113").Value = 1499656.2
$ws.Range("I113").Value = 2299186
$ws.Range("J113").Value = 538
$ws.Range("K113").Value = 6897558
$ws.Range("L113").Value = 1614
$ws.Range("M113").Value = -6895388
$ws.Range("N113").Value = -5954
# Row 129
$ws.Range("H129").Value = 1871.119
$ws.Range("I129").Value = 1704.7858
$ws.Range("J129").Value = 1954.2858
$ws.Range("K129").Value = 5114.357400000001
$ws.Range("L129").Value = 5862.857400000001
$ws.Range("M129").Value = -114.3574000000008
$ws.Range("N129").Value = -15862.8574
# Row 131
$ws.Range("H131").Value = 949.8261
$ws.Range("J131").Value = 949.8261
$ws.Range("L131").Value = 2849.4783
$ws.Range("N131").Value = -12929.4783

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2944457.8
$ws.Range("I122").Value = 7693626.5
$ws.Range("J122").Value = 4496.048
$ws.Range("K122").Value = 23080879.5
$ws.Range("L122").Value = 13488.144
$ws.Range("M122").Value = -23078429.5
$ws.Range("N122").Value = -18388.144
# Row 126
$ws.Range("H126").Value = 17285208
$ws.Range("I126").Value = 16667852
$ws.Range("J126").Value = 18519916
$ws.Range("K126").Value = 50003556
$ws.Range("L126").Value = 55559748
$ws.Range("M126").Value = -50001086
$ws.Range("N126").Value = -55564688

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2500
$ws.Range("I7").Value = 2500
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2388
$ws.Range("N7").ClearContents()
# Row 40
$ws.Range("H40").Value = 3489.5862
$ws.Range("I40").Value = 3061.3333
$ws.Range("J40").Value = 4613.75
$ws.Range("K40").Value = 3061.3333
$ws.Range("L40").Value = 4613.75
$ws.Range("M40").Value = -2925.3333
$ws.Range("N40").Value = -4885.75
# Row 126
$ws.Range("H126").Value = 2500
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5030
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 1836.49
$ws.Range("I132").Value = 1833.7391
$ws.Range("J132").Value = 1868.125
$ws.Range("K132").Value = 5501.2173
$ws.Range("L132").Value = 5604.375
$ws.Range("M132").Value = -2971.2173
$ws.Range("N132").Value = -10664.375

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 74
$ws.Range("H74").Value = 20000
$ws.Range("J74").Value = 20000
$ws.Range("L74").Value = 20000
$ws.Range("N74").Value = -21872
# Row 77
$ws.Range("H77").Value = 20000
$ws.Range("J77").Value = 20000
$ws.Range("L77").Value = 60000
$ws.Range("N77").Value = -69360
# Row 126
$ws.Range("H126").Value = 2510.9487
$ws.Range("I126").Value = 2498.0688
$ws.Range("J126").Value = 2548.3
$ws.Range("K126").Value = 7494.2064
$ws.Range("L126").Value = 7644.900000000001
$ws.Range("M126").Value = -5024.2064
$ws.Range("N126").Value = -12584.9
# Row 132
$ws.Range("H132").Value = 2025.4584
$ws.Range("I132").Value = 1374.3422
$ws.Range("J132").Value = 4499.7
$ws.Range("K132").Value = 4123.0266
$ws.Range("L132").Value = 13499.1
$ws.Range("M132").Value = -1593.0266
$ws.Range("N132").Value = -18559.1
